$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 5) matching the existing data rows.
$ws.Range("A5").Value = 42602.016458333332
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"

$ws.Range("B5").Value = "Named"

$ws.Range("C5").Value = 5933
$ws.Range("D5").Value = 3076
$ws.Range("E5").Value = 194
$ws.Range("F5").Value = 42
$ws.Range("G5").Value = 15
$ws.Range("H5").Value = 73
$ws.Range("I5").Value = 26
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 0
